$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1436.3182
$ws.Range("I15").Value = 1436.3182
$ws.Range("K15").Value = 4308.9546
$ws.Range("M15").Value = -4139.9546

$ws.Range("H33").Value = 316.68182
$ws.Range("I33").Value = 296.78946
$ws.Range("K33").Value = 296.78946
$ws.Range("M33").Value = -67.78946000000002

$ws.Range("H113").Value = 4007.0833
$ws.Range("I113").Value = 4007.7273
$ws.Range("K113").Value = 4007.7273
$ws.Range("M113").Value = -753.7273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5085.7144
$ws.Range("I32").Value = 4470.231
$ws.Range("K32").Value = 4470.231
$ws.Range("M32").Value = -4183.231

$ws.Range("H45").Value = 2165.348
$ws.Range("I45").Value = 1990.25
$ws.Range("J45").Value = 3332.6667
$ws.Range("K45").Value = 1990.25
$ws.Range("L45").Value = 3332.6667
$ws.Range("M45").Value = -1613.25
$ws.Range("N45").Value = -4086.6667

$ws.Range("H61").Value = 1646.6428
$ws.Range("I61").Value = 1496.0834
$ws.Range("J61").Value = 2550
$ws.Range("K61").Value = 1496.0834
$ws.Range("L61").Value = 2550
$ws.Range("M61").Value = -1284.0834
$ws.Range("N61").Value = -2974

$ws.Range("H74").Value = 4426.636
$ws.Range("I74").Value = 3959.8333
$ws.Range("K74").Value = 3959.8333
$ws.Range("M74").Value = -3085.8333

$ws.Range("H77").Value = 4426.636
$ws.Range("I77").Value = 3959.8333
$ws.Range("K77").Value = 19799.1665
$ws.Range("M77").Value = -15431.1665

$ws.Range("H102").Value = 3186.8
$ws.Range("I102").Value = 3108.5
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 3108.5
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -1486.5
$ws.Range("N102").Value = -6744

$ws.Range("H132").Value = 6631.343
$ws.Range("I132").Value = 5654.0415
$ws.Range("K132").Value = 16962.1245
$ws.Range("M132").Value = -14432.1245

$ws.Range("H136").Value = 1646.6428
$ws.Range("I136").Value = 1496.0834
$ws.Range("J136").Value = 2550
$ws.Range("K136").Value = 4488.2502
$ws.Range("L136").Value = 7650
$ws.Range("M136").Value = -1938.2502
$ws.Range("N136").Value = -12750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13738.083
$ws.Range("I86").Value = 23955.4
$ws.Range("J86").Value = 6440
$ws.Range("K86").Value = 23955.4
$ws.Range("L86").Value = 6440
$ws.Range("M86").Value = -22832.4
$ws.Range("N86").Value = -8686

$ws.Range("H89").Value = 13738.083
$ws.Range("I89").Value = 23955.4
$ws.Range("J89").Value = 6440
$ws.Range("K89").Value = 119777
$ws.Range("L89").Value = 32200
$ws.Range("M89").Value = -114161
$ws.Range("N89").Value = -43432

$ws.Range("H99").Value = 1157.5834
$ws.Range("I99").Value = 739.1
$ws.Range("K99").Value = 739.1
$ws.Range("M99").Value = 758.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H31").Value = 9127.5
$ws.Range("I31").Value = 13314.667
$ws.Range("J31").Value = 2846.75
$ws.Range("K31").Value = 13314.667
$ws.Range("L31").Value = 2846.75
$ws.Range("M31").Value = -13019.667
$ws.Range("N31").Value = -3436.75

$ws.Range("H34").Value = 9127.5
$ws.Range("I34").Value = 13314.667
$ws.Range("J34").Value = 2846.75
$ws.Range("K34").Value = 13314.667
$ws.Range("L34").Value = 2846.75
$ws.Range("M34").Value = -13112.667
$ws.Range("N34").Value = -3250.75

$ws.Range("H122").Value = 2111.182
$ws.Range("I122").Value = 2217.8572
$ws.Range("K122").Value = 6653.571599999999
$ws.Range("M122").Value = -4203.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2250
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 9000
$ws.Range("M86").Value = -7814

$ws.Range("H89").Value = 2250
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 27000
$ws.Range("M89").Value = -21072

$ws.Range("H106").Value = 17312.375
$ws.Range("J106").Value = 19416.666
$ws.Range("L106").Value = 58249.99800000001
$ws.Range("N106").Value = -60141.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1651.84
$ws.Range("I102").Value = 1432.85
$ws.Range("J102").Value = 2527.8
$ws.Range("K102").Value = 1432.85
$ws.Range("L102").Value = 2527.8
$ws.Range("M102").Value = 189.1500000000001
$ws.Range("N102").Value = -5771.8

$ws.Range("H113").Value = 5903.6665
$ws.Range("I113").Value = 4100
$ws.Range("K113").Value = 4100
$ws.Range("M113").Value = -1930

$ws.Range("H122").Value = 5007
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3384.7334
$ws.Range("I46").Value = 3415.1724
$ws.Range("K46").Value = 3415.1724
$ws.Range("M46").Value = -3227.1724

$ws.Range("H68").Value = 2739.8
$ws.Range("I68").Value = 3233
$ws.Range("K68").Value = 3233
$ws.Range("M68").Value = -2484

$ws.Range("H71").Value = 2739.8
$ws.Range("I71").Value = 3233
$ws.Range("K71").Value = 16165
$ws.Range("M71").Value = -12421

$ws.Range("H136").Value = 2183.2222
$ws.Range("I136").Value = 2124.1428
$ws.Range("J136").Value = 2390
$ws.Range("K136").Value = 6372.428400000001
$ws.Range("L136").Value = 7170
$ws.Range("M136").Value = -3822.428400000001
$ws.Range("N136").Value = -12270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H81").Value = 2421.7778
$ws.Range("J81").Value = 4000.3333
$ws.Range("L81").Value = 8000.6666
$ws.Range("N81").Value = -10122.6666

$ws.Range("H84").Value = 2421.7778
$ws.Range("J84").Value = 4000.3333
$ws.Range("L84").Value = 40003.333
$ws.Range("N84").Value = -50611.333

$ws.Range("H132").Value = 2161.25
$ws.Range("I132").Value = 2161.25
$ws.Range("K132").Value = 6483.75
$ws.Range("M132").Value = -3953.75
